# Update "想去人数" (number of people interested) counts across sheets.
# The same underlying events appear on both the "展览" sheet and the
# "全部类型" sheet, so each value needs to be bumped in both places.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 72
$ws1.Range("F3").Value = 11773
$ws1.Range("F5").Value = 343
$ws1.Range("F7").Value = 11710
$ws1.Range("F12").Value = 1768
$ws1.Range("F13").Value = 5798
$ws1.Range("F14").Value = 120
$ws1.Range("F15").Value = 3525

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 72
$ws4.Range("F5").Value = 11773
$ws4.Range("F7").Value = 343
$ws4.Range("F9").Value = 11710
$ws4.Range("F14").Value = 1768
$ws4.Range("F16").Value = 5798
$ws4.Range("F17").Value = 120
$ws4.Range("F18").Value = 3525
